# Update the source data on the "DR42-level-of-trust-fairness-an" sheet
# with corrected figures (re-run analysis fixing prior data errors).
# The workbook's pivot table ("PivotTable1" on Sheet1) is sourced from
# this range and recalculates automatically from the live data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DR42-level-of-trust-fairness-an")

# 2002 - Fairness
$ws.Range("N2").Value = 5.5454545454545396
$ws.Range("O2").Value = 2052
$ws.Range("P2").Value = 5.4538150077266403
$ws.Range("Q2").Value = 5.6370940831824496
$ws.Range("R2").Value = 0.04675486618770700181
# 2002 - Helpfulness
$ws.Range("N3").Value = 5.4247311827956999
$ws.Range("O3").Value = 2052
$ws.Range("P3").Value = 5.3347518310843096
$ws.Range("Q3").Value = 5.5147105345070901
$ws.Range("R3").Value = 0.04590783250581179781
# 2002 - Trust
$ws.Range("N4").Value = 5.0473401659346004
$ws.Range("O4").Value = 2052
$ws.Range("P4").Value = 4.9513657917123899
$ws.Range("Q4").Value = 5.1433145401568199
$ws.Range("R4").Value = 0.04896651746031419755
# 2004 - Fairness
$ws.Range("N5").Value = 5.5459373340414198
$ws.Range("O5").Value = 1897
$ws.Range("P5").Value = 5.4532914123852203
$ws.Range("Q5").Value = 5.6385832556976201
$ws.Range("R5").Value = 0.04726832737561219994
# 2004 - Helpfulness
$ws.Range("N6").Value = 5.6663135593220399
$ws.Range("O6").Value = 1897
$ws.Range("P6").Value = 5.5748344607409903
$ws.Range("Q6").Value = 5.7577926579030798
$ws.Range("R6").Value = 0.04667300948012410161
# 2004 - Trust
$ws.Range("N7").Value = 5.1362916006339203
$ws.Range("O7").Value = 1897
$ws.Range("P7").Value = 5.0381822965655099
$ws.Range("Q7").Value = 5.2344009047023201
$ws.Range("R7").Value = 0.05005576738183910301
# 2006 - Fairness
$ws.Range("N8").Value = 5.7802013422818801
$ws.Range("O8").Value = 2394
$ws.Range("P8").Value = 5.6976388389597199
$ws.Range("Q8").Value = 5.8627638456040296
$ws.Range("R8").Value = 0.04212372618477280145
# 2006 - Helpfulness
$ws.Range("N9").Value = 5.6838574423480104
$ws.Range("O9").Value = 2394
$ws.Range("P9").Value = 5.6025661233647099
$ws.Range("Q9").Value = 5.7651487613313099
$ws.Range("R9").Value = 0.04147516274658050173
# 2006 - Trust
$ws.Range("N10").Value = 5.4035234899328799
$ws.Range("O10").Value = 2394
$ws.Range("P10").Value = 5.31647088441989
$ws.Range("Q10").Value = 5.4905760954458698
$ws.Range("R10").Value = 0.04441459464948469688
# 2008 - Fairness
$ws.Range("N11").Value = 5.6821870995301103
$ws.Range("O11").Value = 2352
$ws.Range("P11").Value = 5.5965545601090598
$ws.Range("Q11").Value = 5.7678196389511696
$ws.Range("R11").Value = 0.04369007113319280100
# 2008 - Helpfulness
$ws.Range("N12").Value = 5.62079250106519
$ws.Range("O12").Value = 2352
$ws.Range("P12").Value = 5.5368182596788103
$ws.Range("Q12").Value = 5.7047667424515698
$ws.Range("R12").Value = 0.04284400070733460075
# 2008 - Trust
$ws.Range("N13").Value = 5.2687393526405497
$ws.Range("O13").Value = 2352
$ws.Range("P13").Value = 5.1779302226919004
$ws.Range("Q13").Value = 5.3595484825891901
$ws.Range("R13").Value = 0.04633118874930779868
# 2010 - Fairness
$ws.Range("N14").Value = 5.6193656093489102
$ws.Range("O14").Value = 2422
$ws.Range("P14").Value = 5.5378847815901899
$ws.Range("Q14").Value = 5.7008464371076402
$ws.Range("R14").Value = 0.04157185089730929678
# 2010 - Helpfulness
$ws.Range("N15").Value = 5.7862785862785904
$ws.Range("O15").Value = 2422
$ws.Range("P15").Value = 5.7094247806846399
$ws.Range("Q15").Value = 5.86313239187254
$ws.Range("R15").Value = 0.03921112530303520055
# 2010 - Trust
$ws.Range("N16").Value = 5.3556935817805398
$ws.Range("O16").Value = 2422
$ws.Range("P16").Value = 5.2694984210147604
$ws.Range("Q16").Value = 5.44188874254632
$ws.Range("R16").Value = 0.04397712283968469721
# 2012 - Fairness
$ws.Range("N17").Value = 5.7861552028218703
$ws.Range("O17").Value = 2286
$ws.Range("P17").Value = 5.70667760163722
$ws.Range("Q17").Value = 5.8656328040065198
$ws.Range("R17").Value = 0.04054979652278020014
# 2012 - Helpfulness
$ws.Range("N18").Value = 5.8337730870712399
$ws.Range("O18").Value = 2286
$ws.Range("P18").Value = 5.7540942154058703
$ws.Range("Q18").Value = 5.9134519587366103
$ws.Range("R18").Value = 0.04065248554355579969
# 2012 - Trust
$ws.Range("N19").Value = 5.3862038664323402
$ws.Range("O19").Value = 2286
$ws.Range("P19").Value = 5.2996685740229799
$ws.Range("Q19").Value = 5.4727391588416898
$ws.Range("R19").Value = 0.04415065939252980215
# 2014 - Fairness
$ws.Range("N20").Value = 5.8399822695035501
$ws.Range("O20").Value = 2264
$ws.Range("P20").Value = 5.7552859176087496
$ws.Range("Q20").Value = 5.92467862139834
$ws.Range("R20").Value = 0.04321242443612160339
# 2014 - Helpfulness
$ws.Range("N21").Value = 5.9203539823008899
$ws.Range("O21").Value = 2264
$ws.Range("P21").Value = 5.8360152588950802
$ws.Range("Q21").Value = 6.0046927057066899
$ws.Range("R21").Value = 0.04302996092132699985
# 2014 - Trust
$ws.Range("N22").Value = 5.3756637168141603
$ws.Range("O22").Value = 2264
$ws.Range("P22").Value = 5.2842870703895501
$ws.Range("Q22").Value = 5.4670403632387696
$ws.Range("R22").Value = 0.04662073797174100254
# 2016 - Fairness
$ws.Range("N23").Value = 5.7407975460122698
$ws.Range("O23").Value = 1959
$ws.Range("P23").Value = 5.6495557880663902
$ws.Range("Q23").Value = 5.8320393039581502
$ws.Range("R23").Value = 0.04655191731932759730
# 2016 - Helpfulness
$ws.Range("N24").Value = 5.8065506653019403
$ws.Range("O24").Value = 1959
$ws.Range("P24").Value = 5.7193733618865803
$ws.Range("Q24").Value = 5.8937279687173101
$ws.Range("R24").Value = 0.04447821602824750253
# 2016 - Trust
$ws.Range("N25").Value = 5.3244762391415401
$ws.Range("O25").Value = 1959
$ws.Range("P25").Value = 5.2266838881319799
$ws.Range("Q25").Value = 5.4222685901511003
$ws.Range("R25").Value = 0.04989405663753210263
# 2018 - Fairness
$ws.Range("N26").Value = 5.7492020063839497
$ws.Range("O26").Value = 2204
$ws.Range("P26").Value = 5.6617764785955096
$ws.Range("Q26").Value = 5.8366275341723899
$ws.Range("R26").Value = 0.04460486111655079938
# 2018 - Helpfulness
$ws.Range("N27").Value = 5.7752502274795301
$ws.Range("O27").Value = 2204
$ws.Range("P27").Value = 5.69047387146838
$ws.Range("Q27").Value = 5.8600265834906802
$ws.Range("R27").Value = 0.04325324286283140085
# 2018 - Trust
$ws.Range("N28").Value = 5.17507958162801
$ws.Range("O28").Value = 2204
$ws.Range("P28").Value = 5.0795520173124897
$ws.Range("Q28").Value = 5.2706071459435302
$ws.Range("R28").Value = 0.04873855322220380093
# 2020 - Fairness
$ws.Range("N29").Value = 5.9128160418483002
$ws.Range("O29").Value = 1149
$ws.Range("P29").Value = 5.7944868591251799
$ws.Range("Q29").Value = 6.0311452245714197
$ws.Range("R29").Value = 0.06037203200159319688
# 2020 - Helpfulness
$ws.Range("N30").Value = 5.93554006968641
$ws.Range("O30").Value = 1149
$ws.Range("P30").Value = 5.8187012691138102
$ws.Range("Q30").Value = 6.0523788702590098
$ws.Range("R30").Value = 0.05961163294520240019
# 2020 - Trust
$ws.Range("N31").Value = 5.4255874673629201
$ws.Range("O31").Value = 1149
$ws.Range("P31").Value = 5.2970507491747902
$ws.Range("Q31").Value = 5.5541241855510597
$ws.Range("R31").Value = 0.06557995825925300171

$wb.RefreshAll()
